$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to Text format so Excel keeps the values as literal
# strings (matching the original inline-string cell type) rather than
# auto-converting them to numbers/currency/percentages.
$cells = @("B2", "C2", "D2", "F2", "G2", "I2")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B2").Value = "$71.77"
$ws.Range("C2").Value = "+0.15(0.21%) 1D"
$ws.Range("D2").Value = "$308.03B"
$ws.Range("F2").Value = "29"
$ws.Range("G2").Value = "12.4"
$ws.Range("I2").Value = "2.85%"
